$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 16.74135869433504
$ws.Range("C2").Value = 10.40645936437756
$ws.Range("D2").Value = 9.655582376328857
$ws.Range("E2").Value = 13.86504579679887
$ws.Range("F2").Value = 30.51856598266494
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 9.830610554194411
$ws.Range("O2").Value = 22.67541293113988
$ws.Range("B3").Value = 16.0337889323745
$ws.Range("C3").Value = 9.769872995893483
$ws.Range("D3").Value = 9.586803596301898
$ws.Range("E3").Value = 13.79737846091927
$ws.Range("F3").Value = 30.61092083064
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 9.839267032918404
$ws.Range("O3").Value = 22.80050043012217
$ws.Range("B4").Value = 15.58366149769168
$ws.Range("C4").Value = 9.356148757909725
$ws.Range("D4").Value = 9.545502157929738
$ws.Range("E4").Value = 13.75838719220569
$ws.Range("F4").Value = 30.67833523616045
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 9.846274117652481
$ws.Range("O4").Value = 22.88484312272922
$ws.Range("B5").Value = 15.3965278838583
$ws.Range("C5").Value = 9.18185437567225
$ws.Range("D5").Value = 9.528917775186825
$ws.Range("E5").Value = 13.7431527262215
$ws.Range("F5").Value = 30.70848496132128
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 9.849554982984539
$ws.Range("O5").Value = 22.92109894140573
$ws.Range("B6").Value = 15.36523844454909
$ws.Range("C6").Value = 9.152570486922825
$ws.Range("D6").Value = 9.526179189183996
$ws.Range("E6").Value = 13.74066293238982
$ws.Range("F6").Value = 30.71365257970405
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 9.850125459237217
$ws.Range("O6").Value = 22.92723277414734
$ws.Range("B7").Value = 15.58115240270041
$ws.Range("C7").Value = 9.353821159798466
$ws.Range("D7").Value = 9.545277481460364
$ws.Range("E7").Value = 13.75817906859627
$ws.Range("F7").Value = 30.67873102398608
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 9.846316642171834
$ws.Range("O7").Value = 22.8853244605959
$ws.Range("B8").Value = 16.50077614483926
$ws.Range("C8").Value = 10.19173382986081
$ws.Range("D8").Value = 9.631682356991785
$ws.Range("E8").Value = 13.841190832577
$ws.Range("F8").Value = 30.54817883362026
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 9.833244147035407
$ws.Range("O8").Value = 22.71697242902223
$ws.Range("B9").Value = 18.17079458782194
$ws.Range("C9").Value = 11.65199651999141
$ws.Range("D9").Value = 9.807920635204955
$ws.Range("E9").Value = 14.02370279445567
$ws.Range("F9").Value = 30.37774723420653
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 9.821032360781007
$ws.Range("O9").Value = 22.44709805205695
$ws.Range("B10").Value = 19.30636051691138
$ws.Range("C10").Value = 12.61204613125867
$ws.Range("D10").Value = 9.94076171789939
$ws.Range("E10").Value = 14.16901682208175
$ws.Range("F10").Value = 30.30547356949216
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.820236166297818
$ws.Range("O10").Value = 22.28616836538812
$ws.Range("B11").Value = 19.80146181603461
$ws.Range("C11").Value = 13.02417847875567
$ws.Range("D11").Value = 10.00175217461035
$ws.Range("E11").Value = 14.23737480410094
$ws.Range("F11").Value = 30.28422967106768
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.821645231203954
$ws.Range("O11").Value = 22.22119525824375
$ws.Range("B12").Value = 19.98574531011619
$ws.Range("C12").Value = 13.17670081036793
$ws.Range("D12").Value = 10.02491357235215
$ws.Range("E12").Value = 14.26356765386686
$ws.Range("F12").Value = 30.27786736953861
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.8224329212455
$ws.Range("O12").Value = 22.19778618970866
$ws.Range("B13").Value = 19.94620052110499
$ws.Range("C13").Value = 13.14401003402629
$ws.Range("D13").Value = 10.01992267020315
$ws.Range("E13").Value = 14.25791316902744
$ws.Range("F13").Value = 30.27916265004497
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.822251987524593
$ws.Range("O13").Value = 22.20277444202605
$ws.Range("B14").Value = 19.81668765893341
$ws.Range("C14").Value = 13.0367976137004
$ws.Range("D14").Value = 10.00365644707196
$ws.Range("E14").Value = 14.23952364845333
$ws.Range("F14").Value = 30.28367247621409
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.821704946136043
$ws.Range("O14").Value = 22.21924537303272
$ws.Range("B15").Value = 19.73693742919219
$ws.Range("C15").Value = 12.97066537185796
$ws.Range("D15").Value = 9.993701019038209
$ws.Range("E15").Value = 14.22829902636576
$ws.Range("F15").Value = 30.2866542143398
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.821402939813607
$ws.Range("O15").Value = 22.22949021188946
$ws.Range("B16").Value = 19.27356233458896
$ws.Range("C16").Value = 12.5846167111861
$ws.Range("D16").Value = 9.936785847967494
$ws.Range("E16").Value = 14.16459338536197
$ws.Range("F16").Value = 30.30709684357646
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 9.820179674855899
$ws.Range("O16").Value = 22.29058113146247
$ws.Range("B17").Value = 18.98371088406165
$ws.Range("C17").Value = 12.34148232660937
$ws.Range("D17").Value = 9.90200314595679
$ws.Range("E17").Value = 14.12607743442113
$ws.Range("F17").Value = 30.32262454306314
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 9.8198825207886
$ws.Range("O17").Value = 22.330175185397
$ws.Range("B18").Value = 18.81498034558319
$ws.Range("C18").Value = 12.19932252460653
$ws.Range("D18").Value = 9.882050701333492
$ws.Range("E18").Value = 14.10413727660715
$ws.Range("F18").Value = 30.33265021616471
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 9.819878368973125
$ws.Range("O18").Value = 22.35372308781221
$ws.Range("B19").Value = 18.75750862612264
$ws.Range("C19").Value = 12.15079225009724
$ws.Range("D19").Value = 9.875304814568905
$ws.Range("E19").Value = 14.09674584057944
$ws.Range("F19").Value = 30.33623241797875
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 9.819905621004274
$ws.Range("O19").Value = 22.36182868655836
$ws.Range("B20").Value = 19.01477555135667
$ws.Range("C20").Value = 12.36760408186225
$ws.Range("D20").Value = 9.905700379962298
$ws.Range("E20").Value = 14.13015557896143
$ws.Range("F20").Value = 30.32085823847426
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 9.819896897863904
$ws.Range("O20").Value = 22.325880099958
$ws.Range("B21").Value = 19.85481640185663
$ws.Range("C21").Value = 13.06838468946346
$ws.Range("D21").Value = 10.0084325691042
$ws.Range("E21").Value = 14.24491689763331
$ws.Range("F21").Value = 30.28230210609044
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.821858734463861
$ws.Range("O21").Value = 22.21437494928537
$ws.Range("B22").Value = 20.38513559736473
$ws.Range("C22").Value = 13.50573263669406
$ws.Range("D22").Value = 10.075949616848
$ws.Range("E22").Value = 14.32170224031975
$ws.Range("F22").Value = 30.2669125168948
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.824621670701003
$ws.Range("O22").Value = 22.14846915645766
$ws.Range("B23").Value = 20.10383777703407
$ws.Range("C23").Value = 13.27420243844966
$ws.Range("D23").Value = 10.03988507414156
$ws.Range("E23").Value = 14.28056311540912
$ws.Range("F23").Value = 30.27422597261662
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.823011781666366
$ws.Range("O23").Value = 22.18300313648146
$ws.Range("B24").Value = 19.00073771978765
$ws.Range("C24").Value = 12.35580184474569
$ws.Range("D24").Value = 9.904028720445474
$ws.Range("E24").Value = 14.12831121568393
$ws.Range("F24").Value = 30.32165336338393
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 9.819889878730969
$ws.Range("O24").Value = 22.32781946435491
$ws.Range("B25").Value = 17.73442768632358
$ws.Range("C25").Value = 11.2767034259503
$ws.Range("D25").Value = 9.759595803215394
$ws.Range("E25").Value = 13.97229656838438
$ws.Range("F25").Value = 30.41460723276341
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 9.822899182246912
$ws.Range("O25").Value = 22.51359136460956
